$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contact entries appended to the Contacts sheet (rows 5-8).
# Columns: A=Name, B=Phone, C=Project, D=Message, E=Date
$rows = @(
    @("rihan",     "9684751485", "2BHK", "Hello i need home",                                            "8/3/2025, 1:07:42 pm"),
    @("Tommy",      "2545487474", "2BHK", "Homeless guy want some home ",                                 "11/3/2025, 2:20:47 pm"),
    @("mohan pai",  "9845444444", "3BHK", "Joe less home",                                                "27/3/2025, 12:12:52 pm"),
    @("Richard ",   "9854747474", "2BHK", "apartment purchase in bc road required we want home to live ", "27/3/2025, 1:09:26 pm")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    # Phone numbers are digit strings that must stay text (like existing rows),
    # not be auto-converted to numbers. Prefix with an apostrophe to force text
    # storage, then reset the style so no quote-prefix formatting lingers on
    # the cell.
    $ws.Cells.Item($r, 2).Value = "'" + $data[1]
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
